$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-01-14 Sunday"; new = "2024-01-15 Monday"},
    @{old = "266÷6="; new = "496÷4="},
    @{old = "402÷7="; new = "762÷8="},
    @{old = "822÷7="; new = "176÷2="},
    @{old = "719÷7="; new = "210÷9="},
    @{old = "507÷6="; new = "338÷2="},
    @{old = "676÷8="; new = "564÷8="},
    @{old = "335÷2="; new = "665÷9="},
    @{old = "958÷7="; new = "746÷9="},
    @{old = "133÷5="; new = "820÷4="},
    @{old = "278÷8="; new = "101÷4="},
    @{old = "920÷4="; new = "263÷7="},
    @{old = "637÷5="; new = "546÷2="},
    @{old = "293÷4="; new = "821÷6="},
    @{old = "616÷7="; new = "688÷9="},
    @{old = "193÷7="; new = "267÷5="},
    @{old = "384÷7="; new = "356÷4="},
    @{old = "804÷2="; new = "932÷9="},
    @{old = "125÷6="; new = "491÷6="},
    @{old = "917÷5="; new = "859÷3="},
    @{old = "757÷5="; new = "267÷3="},
    @{old = "976÷9="; new = "284÷6="},
    @{old = "851÷7="; new = "374÷3="},
    @{old = "608÷3="; new = "251÷9="},
    @{old = "202÷8="; new = "742÷3="},
    @{old = "920÷7="; new = "332÷4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
